# Generate Report for Handoff
# Updates the status/report rows for file "b.md" across the Overview, zh-cn
# and de-de worksheets to reflect that a new handoff package has been
# generated (superseding the out-of-date handback).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdcff39e25cc9677c147578a66a92a57fb8aaef7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdd613febfd2f38ff79735a4f494fb38b3a650c4/e2e/b.md."

# --- Overview sheet: row 3 corresponds to b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 16:47:14"

# Column P (16) should render with stored OOXML width="40". The COM
# ColumnWidth property is offset from the stored width by 5/6, so set it
# accordingly to land exactly on 40 after round-tripping.
$colPWidth = 40 - (5/6)

# --- zh-cn sheet: row 3 corresponds to b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Prefixing with an apostrophe keeps "False" stored as text instead of being
# auto-coerced into a boolean cell; resetting the style back to Normal
# afterwards drops the quote-prefix marker picked up along the way.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-20 16:47:09"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = $colPWidth

# --- de-de sheet: row 3 corresponds to b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-20 16:47:14"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = $colPWidth
